$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring in the same formatting used by the 2019 column (M) for the new
# 2020 column (N) before filling in the values, so styles line up exactly
# with the existing yearly columns.
$ws.Range("M4:M16").Copy() | Out-Null
$ws.Range("N4:N16").PasteSpecial(-4122) | Out-Null

# Year header
$ws.Range("N4").Value = 2020

# Data rows (2020 figures)
$ws.Range("N5").Value = 588.7
$ws.Range("N6").Value = 62.2
$ws.Range("N7").Value = 99.4
$ws.Range("N8").Value = 6.1
$ws.Range("N9").Value = "-"
$ws.Range("N10").Value = 71
$ws.Range("N10").NumberFormat = "0.0"
$ws.Range("N11").Value = 136.3
$ws.Range("N12").Value = 103.3
$ws.Range("N13").Value = 103.2
$ws.Range("N14").Value = 1.8
$ws.Range("N15").Value = "-"
$ws.Range("N16").Value = 5.4

# Match the saved selection from the source workbook
$ws.Range("P15").Select() | Out-Null
